$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.358.29'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '2.108.54'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''345.05'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = '''0.5241'
$ws.Range('E7').Value = '  +2.24%  '
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('D9').Value = '''54.63'
$ws.Range('E9').Value = '  +4.51%  '
$ws.Range('D10').Value = '''0.09439'
$ws.Range('E10').Value = '  +4.02%  '
$ws.Range('D11').Value = '''1.173'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '''25.07'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').Value = '''8.773'
$ws.Range('E13').Value = '  +8.09%  '
$ws.Range('D14').Value = '''6.924'
$ws.Range('E14').Value = '  +2.73%  '
$ws.Range('D15').Value = '2.077.46'
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('D16').Value = '''101.79'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '''0.00001165'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '''1.007'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = '''21.29'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '''0.06722'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = '''6.313'
$ws.Range('E21').Value = '  +2.58%  '
$ws.Range('D22').Value = '''1.006'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '30.391.49'
$ws.Range('D24').Value = '''12.65'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').Value = '''2.310'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = '2.331.42'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').Value = '''22.06'
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').Value = '''163.72'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').Value = '''2.538'
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').Value = '''133.67'
$ws.Range('E30').Value = '  +1.31%  '
$ws.Range('D31').Value = '''1.151'
$ws.Range('E31').Value = '  +2.54%  '
$ws.Range('D32').Value = '''1.742'
$ws.Range('E32').Value = '  +7.73%  '
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').Value = '''6.905'
$ws.Range('E34').Value = '  +15.21%  '
$ws.Range('D35').Value = '''6.272'
$ws.Range('E35').Value = '  +2.81%  '
$ws.Range('D36').Value = '''3.930'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').Value = '''10.51'
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('E38').Value = '  +2.49%  '
$ws.Range('D39').Value = '''0.06826'
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('D40').Value = '''0.7052'
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '''12.59'
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.345'
$ws.Range('E42').Value = '  +5.37%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').Value = '''0.6843'
$ws.Range('E44').Value = '  +2.70%  '
$ws.Range('D45').Value = '''14.56'
$ws.Range('E45').Value = '  +3.54%  '
$ws.Range('D46').Value = '''2.360'
$ws.Range('E46').Value = '  +3.81%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').Value = '''1.358'
$ws.Range('E48').Value = '  +16.07%  '
$ws.Range('D49').Value = '''3.649'
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('D50').Value = '''0.00000000346'
$ws.Range('E50').Value = '  +3.41%  '
$ws.Range('D51').Value = '''1.223'
$ws.Range('E51').Value = '  +0.70%  '

# The values above look like plain decimal numbers, so writing them via
# .Value (like typing into Excel) auto-converts them to numbers unless we
# force text with a leading quote-prefix apostrophe. That prefix tags the
# cell style as text/quote-prefixed, so reset each cell's Style back to
# Normal afterwards (looping since Style on a multi-area Range only hits
# the first area) to keep the cell's formatting as it was originally.
$quotedCells = @(
    'D5',
    'D7',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D14',
    'D16',
    'D17',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D24',
    'D25',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D34',
    'D35',
    'D36',
    'D37',
    'D39',
    'D40',
    'D41',
    'D42',
    'D44',
    'D45',
    'D46',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $quotedCells) {
    $ws.Range($addr).Style = 'Normal'
}
